# Fix Training Data Issue: the "Date" column (BF) held the source file's
# name-derived label ("4-30-2007-08") instead of the actual game date.
# NBA stats for this file correspond to 2008-04-30, so rewrite every
# data row's Date cell accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")

# Assign via a literal-text formula first so Excel does not auto-convert
# the ISO-looking string into a date serial number, then convert the
# formulas to static values in place (copy / paste-special values) so
# the cells keep storing plain text, matching the original inline-string
# cells already used throughout the sheet.
$rng.Formula = "=""2008-04-30"""
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = $false
